$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 10.50827633333333
$ws.Range("N2").Value = 31.524829
$ws.Range("O2").Value = 0.1682660991018133
$ws.Range("P2").Value = 0.1682660991018134
$ws.Range("Q2").Value = 0.823789317634111
$ws.Range("R2").Value = 7.414103858707001
$ws.Range("S2").Value = 0.1682660991018133
$ws.Range("T2").Value = 0.1682660991018134

# Row 3
$ws.Range("O3").Value = 0.4955285863849104
$ws.Range("P3").Value = 0.4955285863849105
$ws.Range("S3").Value = 0.4955285863849104
$ws.Range("T3").Value = 0.4955285863849105

# Row 4
$ws.Range("M4").Value = 6.495209666666667
$ws.Range("N4").Value = 19.485629
$ws.Range("O4").Value = 0.1040059814559238
$ws.Range("P4").Value = 0.1040059814559238
$ws.Range("Q4").Value = 0.5091876316785555
$ws.Range("R4").Value = 4.582688685107
$ws.Range("S4").Value = 0.1040059814559238
$ws.Range("T4").Value = 0.1040059814559238

# Row 5
$ws.Range("M5").Value = 9.909791666666667
$ws.Range("N5").Value = 29.729375
$ws.Range("O5").Value = 0.1586827309986352
$ws.Range("P5").Value = 0.1586827309986352
$ws.Range("Q5").Value = 0.7768715111805555
$ws.Range("R5").Value = 6.991843600625001
$ws.Range("S5").Value = 0.1586827309986352
$ws.Range("T5").Value = 0.1586827309986352

# Row 6
$ws.Range("M6").Value = 4.591137333333333
$ws.Range("N6").Value = 13.773412
$ws.Range("O6").Value = 0.07351660205871713
$ws.Range("P6").Value = 0.07351660205871713
$ws.Range("Q6").Value = 0.3599191504884444
$ws.Range("R6").Value = 3.239272354396
$ws.Range("S6").Value = 0.07351660205871713
$ws.Range("T6").Value = 0.07351660205871713
